$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.024.71"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "2.231.79"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  -1.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.555"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.99%  "
$ws.Range("E10").Value = "  -5.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0781"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.83%  "
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "2.571.41"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "2.229.64"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.780"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.71%  "
$ws.Range("D18").Value = "43.901.89"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "0.0₃0906"
$ws.Range("E19").Value = "  -6.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("E24").Value = "  -5.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("E26").Value = "  -7.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("E28").Value = "  -2.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "151.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("E32").Value = "  -10.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0755"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.90%  "
$ws.Range("E34").Value = "  -5.65%  "
$ws.Range("E35").Value = "  -4.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.102"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0302"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.80%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.36%  "
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("D44").Value = "1.840.30"
$ws.Range("E44").Value = "  +4.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.72%  "
$ws.Range("E46").Value = "  -7.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "14.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.97%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "67.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.33%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "73.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "94.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.97%  "
$ws.Range("D51").Value = "2.452.74"
$ws.Range("E51").Value = "  -0.88%  "
